$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.702.16"
$ws.Range("E2").Value = "  +2.50%  "

$ws.Range("D3").Value = "3.556.43"
$ws.Range("E3").Value = "  +1.48%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "581.58"
$ws.Range("E5").Value = "  +1.99%  "

$ws.Range("D6").Value = "187.45"
$ws.Range("E6").Value = "  +1.71%  "

$ws.Range("E7").Value = "  +2.34%  "

$ws.Range("D8").Value = "3.545.11"
$ws.Range("E8").Value = "  +1.48%  "

$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("D10").Value = "0.222"
$ws.Range("E10").Value = "  +18.75%  "

$ws.Range("D11").Value = "0.649"
$ws.Range("E11").Value = "  +0.23%  "

$ws.Range("D12").Value = "54.78"
$ws.Range("E12").Value = "  +1.41%  "

$ws.Range("E13").Value = "  +6.48%  "

$ws.Range("D14").Value = "9.52"
$ws.Range("E14").Value = "  +1.16%  "

$ws.Range("D15").Value = "4.119.26"
$ws.Range("E15").Value = "  +1.25%  "

$ws.Range("D16").Value = "70.718.63"
$ws.Range("E16").Value = "  +2.53%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "'12.80"
$ws.Range("E17").Value = "  +4.74%  "

$ws.Range("E18").Value = "  -0.99%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.552.67"
$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("D20").Value = "575.05"
$ws.Range("E20").Value = "  +6.64%  "

$ws.Range("E21").Value = "  +0.85%  "

$ws.Range("E22").Value = "  -0.72%  "

$ws.Range("D23").Value = "17.72"
$ws.Range("E23").Value = "  -4.54%  "

$ws.Range("E24").Value = "  +3.40%  "

$ws.Range("E25").Value = "  -1.86%  "

$ws.Range("D26").Value = "93.98"
$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("E27").Value = "  +4.32%  "

$ws.Range("E28").Value = "  +2.01%  "

$ws.Range("D29").Value = "9.29"
$ws.Range("E29").Value = "  +2.02%  "

$ws.Range("D30").Value = "32.53"
$ws.Range("E30").Value = "  +2.43%  "

$ws.Range("E31").Value = "  -0.37%  "

$ws.Range("D32").Value = "12.32"
$ws.Range("E32").Value = "  -1.67%  "

$ws.Range("E33").Value = "  +2.06%  "

$ws.Range("B34").Value = "dogwifhat"
$ws.Range("C34").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D34").Value = "'3.80"
$ws.Range("E34").Value = "  +21.00%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "63.05"
$ws.Range("E35").Value = "  -2.75%  "

$ws.Range("D36").Value = "3.35"
$ws.Range("E36").Value = "  +11.45%  "

$ws.Range("D37").Value = "541.46"
$ws.Range("E37").Value = "  -3.51%  "

$ws.Range("E38").Value = "  +4.10%  "

$ws.Range("D39").Value = "38.14"
$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("D40").Value = "0.0₃0805"
$ws.Range("E40").Value = "  +5.33%  "

$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").Value = "3.575.87"
$ws.Range("E42").Value = "  +11.26%  "

$ws.Range("D43").Value = "0.139"
$ws.Range("E43").Value = "  +3.89%  "

$ws.Range("E44").Value = "  +3.24%  "

$ws.Range("D45").Value = "0.0472"
$ws.Range("E45").Value = "  +7.47%  "

$ws.Range("E46").Value = "  -0.77%  "

$ws.Range("E47").Value = "  -1.23%  "

$ws.Range("D48").Value = "9.33"
$ws.Range("E48").Value = "  +4.22%  "

$ws.Range("E49").Value = "  +3.12%  "

$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("B51").Value = "OceanProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D51").Value = "1.47"
$ws.Range("E51").Value = "  +8.54%  "

